$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a string value to a cell while guaranteeing it is stored as
# text (not auto-converted to a number), and without altering the cell's
# existing style. We do this by temporarily remembering the style, writing
# the value with a leading apostrophe (Excel's "force text" prefix), and
# then restoring the original style.
function Set-TextValue($range, [string]$val) {
    $savedStyle = $range.Style
    $range.Value = "'" + $val
    $range.Style = $savedStyle
}

$ws.Range("D2").Value = '59.175.50'
$ws.Range("E2").Value = '  +4.33%  '

$ws.Range("D3").Value = '3.343.65'
$ws.Range("E3").Value = '  +2.77%  '

$ws.Range("E4").Value = '  -0.05%  '

Set-TextValue $ws.Range("D5") '412.96'
$ws.Range("E5").Value = '  +3.71%  '

Set-TextValue $ws.Range("D6") '111.69'
$ws.Range("E6").Value = '  +0.35%  '

$ws.Range("E7").Value = '  +4.70%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("E9").Value = '  +2.21%  '

Set-TextValue $ws.Range("D10") '40.20'
$ws.Range("E10").Value = '  +1.84%  '

Set-TextValue $ws.Range("D11") '0.0982'
$ws.Range("E11").Value = '  +2.73%  '

$ws.Range("E12").Value = '  +1.44%  '

$ws.Range("D13").Value = '3.871.68'
$ws.Range("E13").Value = '  +2.96%  '

Set-TextValue $ws.Range("D14") '8.48'
$ws.Range("E14").Value = '  +4.47%  '

Set-TextValue $ws.Range("D15") '19.42'
$ws.Range("E15").Value = '  +0.63%  '

$ws.Range("D16").Value = '3.346.82'
$ws.Range("E16").Value = '  +2.52%  '

$ws.Range("E17").Value = '  -0.45%  '

$ws.Range("D18").Value = '58.957.14'
$ws.Range("E18").Value = '  +4.24%  '

Set-TextValue $ws.Range("D19") '10.87'
$ws.Range("E19").Value = '  -2.70%  '

$ws.Range("E20").Value = '  +0.90%  '

$ws.Range("E21").Value = '  +5.27%  '

Set-TextValue $ws.Range("D22") '13.13'
$ws.Range("E22").Value = '  +0.76%  '

Set-TextValue $ws.Range("D23") '303.96'
$ws.Range("E23").Value = '  +1.28%  '

Set-TextValue $ws.Range("D24") '75.59'
$ws.Range("E24").Value = '  +0.00%  '

$ws.Range("E25").Value = '  -0.10%  '

Set-TextValue $ws.Range("D26") '28.69'
$ws.Range("E26").Value = '  +1.13%  '

$ws.Range("E27").Value = '  +2.70%  '

Set-TextValue $ws.Range("D28") '7.99'
$ws.Range("E28").Value = '  -1.80%  '

Set-TextValue $ws.Range("D29") '7.48'
$ws.Range("E29").Value = '  +1.45%  '

$ws.Range("E30").Value = '  +0.26%  '

Set-TextValue $ws.Range("D31") '0.116'
$ws.Range("E31").Value = '  +3.88%  '

Set-TextValue $ws.Range("D32") '11.66'
$ws.Range("E32").Value = '  +5.05%  '

Set-TextValue $ws.Range("D33") '1.00'
$ws.Range("E33").Value = '  +0.01%  '

Set-TextValue $ws.Range("D34") '40.06'
$ws.Range("E34").Value = '  +9.06%  '

Set-TextValue $ws.Range("D35") '0.0527'
$ws.Range("E35").Value = '  +8.38%  '

$ws.Range("E36").Value = '  +0.63%  '

Set-TextValue $ws.Range("D37") '52.10'
$ws.Range("E37").Value = '  +0.71%  '

$ws.Range("E38").Value = '  +0.46%  '

Set-TextValue $ws.Range("D39") '0.999'
$ws.Range("E39").Value = '  -0.02%  '

Set-TextValue $ws.Range("D40") '3.47'
$ws.Range("E40").Value = '  -2.14%  '

Set-TextValue $ws.Range("D41") '138.00'
$ws.Range("E41").Value = '  +2.66%  '

$ws.Range("E42").Value = '  +1.87%  '

Set-TextValue $ws.Range("D43") '4.05'
$ws.Range("E43").Value = '  +1.56%  '

$ws.Range("E44").Value = '  -0.55%  '

Set-TextValue $ws.Range("D45") '17.04'
$ws.Range("E45").Value = '  -3.37%  '

Set-TextValue $ws.Range("D46") '0.280'
$ws.Range("E46").Value = '  -1.96%  '

Set-TextValue $ws.Range("D47") '2.26'
$ws.Range("E47").Value = '  +8.37%  '

Set-TextValue $ws.Range("D48") '22.58'
$ws.Range("E48").Value = '  +1.35%  '

$ws.Range("D49").Value = '2.206.64'
$ws.Range("E49").Value = '  +2.67%  '

Set-TextValue $ws.Range("D50") '2.40'
$ws.Range("E50").Value = '  -0.64%  '

Set-TextValue $ws.Range("D51") '1.93'
$ws.Range("E51").Value = '  -12.17%  '
